# Fix the misspelled base name "Incerlick" -> "Incirlick" in the JFACC
# paragraph. Word's spell-checker had wrapped the misspelling in
# <w:proofErr w:type="spellStart"/>...<w:proofErr w:type="spellEnd"/>
# markers; now that the word is spelled correctly those markers must be
# removed as well (otherwise Word would still flag it as a checked-but-
# once-bad word). We rebuild the whole containing paragraph's run content
# via Range.InsertXML so the proof-of-error bookmarks bracketing the old
# run disappear along with the fix, while every other run/paragraph in
# the document (and the rest of this paragraph's runs/formatting) is left
# completely untouched.

$d = $word.ActiveDocument

# Locate the misspelled word.
$found = $d.Content
$found.Find.Execute("Incerlick", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

if (-not $found.Find.Found) {
    throw "Could not find 'Incerlick' in the document"
}

# Expand to the whole enclosing paragraph, then back off the trailing
# paragraph mark so we only replace the paragraph's run content (not the
# paragraph mark itself, which would risk merging it with the next one).
$paraRange = $d.Range($found.Start, $found.End)
$paraRange.Expand(4) | Out-Null        # wdParagraph = 4
$paraRange.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1, back off the pilcrow

$newParagraphXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r w:rsidRPr="008336C6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Joint Force Air Component Command </w:t></w:r>
            <w:r w:rsidR="008336C6" w:rsidRPr="008336C6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>consists</w:t></w:r>
            <w:r w:rsidRPr="008336C6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of the 132nd Virtual Wing split on </w:t></w:r>
            <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Incirlick</w:t></w:r>
            <w:r w:rsidRPr="008336C6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and Ramat David Airbase. E-3 and Tanker squadrons have also arrived in theater and are getting ready to support operations.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$paraRange.InsertXML($newParagraphXml)
